$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.15480000000001
$ws.Range("A14").Value = -21.86869999999999
$ws.Range("A21").Value = -20.01839999999998
$ws.Range("B22").Value = 9.729799999999997
$ws.Range("A23").Value = -19.95379999999998
$ws.Range("B24").Value = 5.959700000000001
$ws.Range("A25").Value = -21.8275
$ws.Range("A26").Value = -21.05459999999997
$ws.Range("B28").Value = 5.748100000000002
$ws.Range("A29").Value = -21.02869999999997
$ws.Range("B36").Value = 9.322500000000009
$ws.Range("B45").Value = 5.259800000000005
$ws.Range("B48").Value = 6.108199999999999
$ws.Range("B49").Value = 5.719099999999995
$ws.Range("B52").Value = 5.2044
$ws.Range("A53").Value = -21.91579999999999
$ws.Range("B53").Value = 5.433299999999997
$ws.Range("B54").Value = 4.833700000000002
$ws.Range("A57").Value = -22.41540000000002
$ws.Range("A59").Value = -21.97209999999999
$ws.Range("A69").Value = -21.6585
$ws.Range("B70").Value = 7.107300000000001
$ws.Range("A79").Value = -20.20090000000001
$ws.Range("A83").Value = -21.824
$ws.Range("B86").Value = 5.312200000000002
$ws.Range("B87").Value = 5.709199999999999
$ws.Range("B89").Value = 4.639999999999995
$ws.Range("A91").Value = -20.55419999999999
$ws.Range("A93").Value = -21.31630000000001
$ws.Range("B101").Value = 4.291799999999999
$ws.Range("A103").Value = -21.7573

$wb.Save()
